# Freelance Timesheet - log two new entries for the Google Calendar
# availability-integration work, then restore the sheet's scroll/zoom/
# selection state the way the author left the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: Sat 2023-02-25, 15:50 - 16:50 -----------------------------
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("A18").Value = 44982
$ws.Range("B18").NumberFormat = "h:mm"
$ws.Range("B18").Value = 0.65972222222222221
$ws.Range("C18").NumberFormat = "h:mm"
$ws.Range("C18").Value = 0.70138888888888884
$ws.Range("D18").Value = "Began work on Google Calendar integration for availability scheduling."
$ws.Range("E18").Value = 1
$ws.Rows.Item(18).RowHeight = 28

# --- Row 19: Sun 2023-02-26, 12:54 - 13:54 -----------------------------
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("A19").Value = 44983
$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("B19").Value = 0.53749999999999998
$ws.Range("C19").NumberFormat = "h:mm"
$ws.Range("C19").Value = 0.57916666666666672
$ws.Range("D19").Value = "Can now update availability in DB based on upcoming week's calendar events. Can also populate upcoming week in calendar with availability stored in DB. Hard-coded for one tutor."
$ws.Range("E19").Value = 1
$ws.Rows.Item(19).RowHeight = 70

# The Total Hours / Total Amount Due formulas (E32, E33) recompute on
# their own since they already sum the E8:E31 range.

# --- Restore the window's scroll position, zoom, and active cell ------
$window = $excel.ActiveWindow
$window.ScrollRow = 5
$window.ScrollColumn = 1
$window.Zoom = 120
$null = $ws.Range("F19").Select()
